# Re-applies the "良牙夏典" event (inserted 2024-07-27, before the existing
# 2024-08-03 "蔚蓝档案only" row) plus a handful of attendee-count bumps to
# the "展览" and "全部类型" sheets, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

function Insert-LiangyaRow($ws, $insertRow) {
    # The row currently at $insertRow (e.g. "蔚蓝档案only") keeps its own
    # sequence number, which the brand-new row takes over; the displaced
    # row is renumbered one higher.
    $seq = [int]($ws.Cells.Item($insertRow, 1).Value())

    # Shift rows at/after $insertRow down by one.
    $ws.Rows.Item($insertRow).Insert()

    $ws.Cells.Item($insertRow + 1, 1).Value = $seq + 1

    # Give the new A cell the same (bordered/centered) look as its neighbors.
    $ws.Range($ws.Cells.Item($insertRow - 1, 1), $ws.Cells.Item($insertRow - 1, 1)).Copy()
    $ws.Range($ws.Cells.Item($insertRow, 1), $ws.Cells.Item($insertRow, 1)).PasteSpecial(-4122)
    $ws.Cells.Item($insertRow, 1).Value = $seq

    # Force column B to stay plain text (otherwise "2024-07-27" gets read
    # back as a date serial), then restore the unstyled look other date
    # cells in the column use.
    $ws.Cells.Item($insertRow, 2).NumberFormat = "@"
    $ws.Cells.Item($insertRow, 2).Value = "2024-07-27"
    $ws.Range($ws.Cells.Item($insertRow - 1, 2), $ws.Cells.Item($insertRow - 1, 2)).Copy()
    $ws.Range($ws.Cells.Item($insertRow, 2), $ws.Cells.Item($insertRow, 2)).PasteSpecial(-4122)

    $ws.Cells.Item($insertRow, 3).Value = "南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）"
    $ws.Cells.Item($insertRow, 4).Value = "民族大道106号 南宁国际会展中心"
    $ws.Cells.Item($insertRow, 5).Value = "2024.07.27 09:30-07.28 17:30"
    $ws.Cells.Item($insertRow, 6).Value = 7
    $ws.Cells.Item($insertRow, 7).Value = 55
    $ws.Cells.Item($insertRow, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85264"
    $ws.Cells.Item($insertRow, 9).Value = "//i0.hdslb.com/bfs/openplatform/202405/dZVcS7eE1715155418142.jpeg"

    $excel.CutCopyMode = 0
}

# --- Sheet "展览" (exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 1948
$wsExhibit.Range("F9").Value = 996
$wsExhibit.Range("F10").Value = 188
Insert-LiangyaRow $wsExhibit 11

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 1948
$wsAll.Range("F10").Value = 996
$wsAll.Range("F11").Value = 188
Insert-LiangyaRow $wsAll 12
